$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.840.05"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -0.70%  '

$ws.Range("D3").Value = "'2.749.13"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -0.78%  '

$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").Value = "'574.20"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -1.17%  '

$ws.Range("D6").Value = "'157.65"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +1.13%  '

$ws.Range("E7").Value = '  -0.06%  '

$ws.Range("E8").Value = '  -1.47%  '

$ws.Range("E9").Value = '  -3.20%  '

$ws.Range("E10").Value = '  -0.06%  '

$ws.Range("D11").Value = "'0.382"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -2.54%  '

$ws.Range("D12").Value = "'5.61"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -17.00%  '

$ws.Range("D13").Value = "'3.232.51"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -0.88%  '

$ws.Range("D14").Value = "'26.44"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -1.49%  '

$ws.Range("D15").Value = "'63.520.30"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -1.04%  '

$ws.Range("D16").Value = "'0.0000150"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -2.33%  '

$ws.Range("D17").Value = "'2.748.80"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -0.73%  '

$ws.Range("D18").Value = "'12.12"
$ws.Range("D18").ClearFormats()

$ws.Range("D19").Value = "'4.80"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -2.16%  '

$ws.Range("D20").Value = "'354.47"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -2.22%  '

$ws.Range("D21").Value = "'6.72"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -4.12%  '

$ws.Range("D22").Value = "'0.537"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.80%  '

$ws.Range("E23").Value = '  -0.37%  '

$ws.Range("D24").Value = "'65.08"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -2.43%  '

$ws.Range("E25").Value = '  -1.19%  '

$ws.Range("E26").Value = '  +0.05%  '

$ws.Range("D27").Value = "'8.37"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -2.25%  '

$ws.Range("E28").Value = '  -0.15%  '

$ws.Range("E29").Value = '  -3.82%  '

$ws.Range("D30").Value = "'6.94"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -2.70%  '

$ws.Range("D31").Value = "'169.03"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -2.50%  '

$ws.Range("D32").Value = "'1.20"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -7.06%  '

$ws.Range("D33").Value = "'20.14"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -2.04%  '

$ws.Range("E34").Value = '  +0.02%  '

$ws.Range("E35").Value = '  -0.55%  '

$ws.Range("E36").Value = '  -0.92%  '

$ws.Range("E37").Value = '  -2.48%  '

$ws.Range("D38").Value = "'0.977"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -3.70%  '

$ws.Range("D39").Value = "'6.16"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +5.18%  '

$ws.Range("D40").Value = "'4.13"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -3.08%  '

$ws.Range("D41").Value = "'325.27"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -5.85%  '

$ws.Range("D42").Value = "'38.79"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -1.38%  '

$ws.Range("D43").Value = "'21.28"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -3.37%  '

$ws.Range("D44").Value = "'0.0584"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -1.82%  '

$ws.Range("D45").Value = "'21.26"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -3.72%  '

$ws.Range("D46").Value = "'0.0253"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -1.51%  '

$ws.Range("D47").Value = "'134.87"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -2.45%  '

$ws.Range("B48").Value = 'Stellar'
$ws.Range("C48").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D48").Value = "'0.101"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -0.90%  '

$ws.Range("B49").Value = 'Mantle'
$ws.Range("C49").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D49").Value = "'0.621"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -4.44%  '

$ws.Range("E50").Value = '  -0.18%  '

$ws.Range("D51").Value = "'11.03"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +0.24%  '
